
$wb = $excel.ActiveWorkbook
$wsNutri = $wb.Worksheets.Item("NutritionalData")
$wsMeasures = $wb.Worksheets.Item("researchMeasures")

# Append sentence to the existing 3/31 diary narrative cell (Z77)
$wsMeasures.Range("Z77").Value = 'Woke up at 1240 am ish, in lots of pain, same as other night, after a bunch of blood and a clot in undies, had to change pad, the lower abdominal and low back ached, drank a bottle of water entirely in two parts, walked around and tried sitting and going to bed but it hurt a lot, sharp pain all over lower midsection. Not much blood passed after the huge heavy flow that woke me. I took a shower and felt a little better by the time done washing hair and body and for the water to go through my system. Braided my hair and was able to go to sleep by 145 am. I woke up by alarm at 530 am and gave Mr. Growly his meds, fed the babies, drank my coffee, and worked on the lab from last week due tomorrow or by the end of the day today. Wasn''t able to determine the unknowns but got the video notes in for descriptions of reactants. I have to compare to the reactions in part 1 observed. Took measurements after 730 am and before eating and after finishing 2nd cup of coffee. Had a reg solid BM reg color as normal after 1st cup of coffee and before my 630 alarm went off. The roommate had just arrived a few minutes before that time. At work ate a chocolate protein shake with 2 bananas with tips cut off and 1 serving protein powder and about 2 cups of almond milk. Then after work, had a bowl of pasta no cheese the last of it and a 4th cup of coffee with 1 tbs organic sugar and 1/4 cup almond milk. Made me too relaxed and had to shake it off, my eyelids felt heavy. Went to clients and got their 15 minutes late due to traffic and started another 15 minutes late. She is super nice though. She will be busy and get back to me on scheduling her next massage. Her cute 2 year old was in the room with us bored because his ipad wasn''t charged and playing with all sorts of stuff like his big blocks and other random parts of toys. He likes to put a bunch of stuff on her to cover her, like pieces of paper, the blankets already in the room, her foam rollers, etc. I ran out of eucalyptus oil and will have to remember to get some more at Sprouts or Marshal''s. Rather be Sprouts but the line to wait is long. I haven''t been there since the quarantine lifted some. Only the Chino Hills one has a line. Not the Corona Sprouts. For Breakfast I had a bowl of pasta and no cheese, because I ate it fast. Was able to send my super nice client yesterday her notes and reciept, but the jotform submissions to download was filled in all blank and emailed her my forwarded email of the jotform receipt that gets sent. I now have a 5 hour family massage with a previous family of two plus another two from a referral. Great people. Sunday I don''t work at ME but will study if I don''t drop the courses, and have a couples massage for 3 hours in Chino Hills at 6 pm. Also a very nice couple. Got home around 8 pm, and still on my rag, spotty, but afraid to drink a glass of wine in case I wake up again like the insides of my belly are kicking me and my low back. Disinfected the supplies, used hot stones and was going to use the massage gun but she didn''t today because was tender in muscles, and got my oily hands on it, so it had to be disinfected. Also, tried out using the hotpacks but in hot water instead, and they didn''t stay warm in the pot of water I brought, and 35 minutes of traffic. I left around 10 minutes till 5 pm and got there at 545 pm. I was supposed to be there at 530 pm, the south 15 is backed up at that time and now with the express lanes, its worse. I am going to try out my hydocollator by plugging it in and using it on my Sunday clients. The other family massage wants me to go up the stairs and its too much to use any add-ons to bring up the stairs too. I finished the soap notes and receipt for today''s client, and emailed admissions at my college to ask them to give me an excused withdraw from my chem course due to the stress and history of cancer, her unfair tests/exams, how heavy they are weighted, being distracted by the camera while taking exams, etc. I have a course I can actually and realistically do well in quiz tomorrow. I want to get great sleep tonight, forget about the lab and the class for now, and study for my other course quiz. It is very interesting, so is chemistry, but the course is designed to filter out for elite super test takers with no full time job, or other factors, and spend way too much time and forget all stress and constraints while interpreting exam questions to recall how to answer them. Ate some quesadillas with the Target Good&Gather mozzparmesan blend, it was delicious, 3 of them with paprika. Shared with Princess. Went to bed around 10 pm.'

# New day row (4/1/2021) on researchMeasures sheet, row 78
$wsMeasures.Range("A78").Value = "Thur"
$wsMeasures.Range("B78").Value = 14
$wsMeasures.Range("C78").Value = DateSerial(2021,4,1)
$wsMeasures.Range("D78").Value = TimeSerial(8,30,0)
$wsMeasures.Range("E78").Value = 72
$wsMeasures.Range("F78").Value = 0
$wsMeasures.Range("G78").Value = 0
$wsMeasures.Range("H78").Value = 0
$wsMeasures.Range("I78").Value = 0
$wsMeasures.Range("J78").Value = TimeSerial(8,30,0)
$wsMeasures.Range("K78").Value = 138.8
$wsMeasures.Range("L78").Formula = "=K78-K77"
$wsMeasures.Range("M78").Formula = "=AB77"
$wsMeasures.Range("N78").Value = 31.5
$wsMeasures.Range("O78").Value = 32.75
$wsMeasures.Range("P78").Value = 10.5
$wsMeasures.Range("Q78").Value = 10.375
$wsMeasures.Range("R78").Value = 19.75
$wsMeasures.Range("S78").Value = 19.75
$wsMeasures.Range("T78").Value = 13
$wsMeasures.Range("U78").Value = 13
$wsMeasures.Range("V78").Value = 17
$wsMeasures.Range("W78").Value = 15
$wsMeasures.Range("X78").Value = 7
$wsMeasures.Range("Y78").Value = 7
$wsMeasures.Range("Z78").Value = 'Woke up at 430 and laid in bed until 10 minutes before 5 am and got up, cleaned pet messes, took babies outside, warmer weather this week, not chilly like previous week. Made their food and Mr. Growly''s meds, one pill he used up of the one he takes every 12 hours. He still has the other 2 meds. I had a cup of coffee, no BM movement until 3rd cup, I drank water and had the 3rd cup of coffee. The other plumber or worker of the landlords came by right before I felt like having a BM, with some movement, but feeling dehydrated. He was rude, because the gate is locked, but I unlocked it and asked him to keep the gate closed for our dogs to not escape, and he didn''t acknowlege the request, and I asked him after unlocking it and he said he would. He''s an old guy in late 60s or 70s even. The landlord has some older workers, I think his bug spray guy is 80 years old. I have class at 9 am and was working on the quiz 2 study guid of ch 5 on sex determining chromosomes and factors. I finished by 830 am and took my measurements, still no BM. I felt somewhat brain fogged or bogged to start off the day, not alert and ready. But thankfully no heavy and painful menstruation in the middle of the night that felt like I was ran over my midsection the last two other nights. Had a BM after entering the data for measurements and updating notes, made another 3 quesadillas before lecture in genetics. Earlier while getting over brain fog or minimal brain clarity less so than normal at the beginning of the day I finished the laundry and prepacked the linen and top cover sets and disinfected hot stones and oils to put in the work van. The hot pads are still drying and were still drying at that time. This house is ca. 1908-1911 Spanish Flu era part of the historical district of Corona, CA, and the workers the landlord has have probably spent a lot of their younger years working on this house. It is odd to have workers that old though. Did the lecture, more of a Q&A, only 4-5 females, some questions, ended early. Found documents to prove my stress related oncology for hodgkins and uploaded those and the form to norco college to remove me from CHE-1A with an EW excused withdrawal because it is stressing me out. Sent that by 10:22 am. The landlord stopped by too to do something with the other side or basement not sure what. He was gone by the time I sent it. Then made flashcards and got ready for quiz. Took the quiz, and feel like I did well on it. Some questions are input types that require actual grading as the other quizes are. I am going to continue the chemistry course, but not invest too much stress into it or my time, if late, then late. I am waiting for the excused withdrawal to be approved, and if not then I have to complete the course so might as well. Made pasta and then logged into course. Thinking about not, and finishing lab, but don''t know yet. I also had a sweet cream nitro cold brew from Starbucks when I went to the bank to get a cashier''s check from the roommate''s share of rent he had me transfer to my checking account. Put that on the porch for the landlord to grab when he is able to. Expecting my decals from build a sign to come in soon. Also, I started out spotty, but when I left to get the cashier''s check and coffee, a bunch of blood fell out and soaked light pad, got on driver seat. Had to use some paper towels after entering the bank and thankfully not have it drip down my leg. I used the paper towels while driving back home and waiting in Starbuck''s drive thru. They put more sweet cream than coffee in my nitro sweet cream. '
$wsMeasures.Range("AA78").Value = '2 bowls of pasta 4/1/2021 recipe
(1494	57.6	16.4	90.2	162.2	34.4	1850)
1/4 cup mozz Stater brand
(80	6	3	6	2	0	180)
1/4 cup mozzParm Good&Gather brand
(100	6	4	8	2	0	280)
3 corn mozzparm blend quesadillas
6 corn tortillas Guerrero
(300	3	0	6	63	6	60)
1/2 cup mozz/parm Good&Gather brand
(200	12	8	16	4	0	560)
Grande Nitro Sweet Cream approximately because it looked mostly cream
(70    5    3.5      1      4     0     20)
=1494+80+100+300+200+70
=57.6+6+6+3+12+5
=16.4+3+4+0+8+3.5
=90.2+6+8+6+16+1
=162.2+2+2+63+4+4
=34.4+0+0+6+0+0
=1850+180+280+60+560+20
'
$wsMeasures.Range("AB78").Formula = "=1494+80+100+300+200+70"
$wsMeasures.Range("AC78").Formula = "=57.6+6+6+3+12+5"
$wsMeasures.Range("AD78").Formula = "=16.4+3+4+0+8+3.5"
$wsMeasures.Range("AE78").Formula = "=90.2+6+8+6+16+1"
$wsMeasures.Range("AF78").Formula = "=162.2+2+2+63+4+4"
$wsMeasures.Range("AG78").Formula = "=34.4+0+0+6+0+0"
$wsMeasures.Range("AH78").Formula = "=1850+180+280+60+560+20"
$wsMeasures.Range("AI78").Formula = "=$AC78/$AB78"
$wsMeasures.Range("AJ78").Formula = "=$AD78/$AB78"
$wsMeasures.Range("AK78").Formula = "=$AE78/$AB78"
$wsMeasures.Range("AL78").Formula = "=$AF78/$AB78"
$wsMeasures.Range("AM78").Formula = "=$AG78/$AB78"
$wsMeasures.Range("AN78").Formula = "=$AH78/$AB78"
$wsMeasures.Range("AO78").Value = 4
$wsMeasures.Range("AP78").Value = 1
$wsMeasures.Range("AQ78").Value = 1
$wsMeasures.Range("AR78").Value = 0
$wsMeasures.Range("AS78").Value = 0
$wsMeasures.Range("AT78").Value = 0
$wsMeasures.Range("AU78").Value = 0
$wsMeasures.Range("AV78").Value = 0
$wsMeasures.Range("AW78").Value = 0
$wsMeasures.Range("AX78").Value = 0
$wsMeasures.Range("AY78").Value = 6.5
$wsMeasures.Range("AZ78").Value = 1
$wsMeasures.Range("BA78").Value = 1
$wsMeasures.Range("BB78").Value = 0
$wsMeasures.Range("BC78").Value = 1
$wsMeasures.Range("BD78").Value = 1
$wsMeasures.Range("BE78").Value = 0
$wsMeasures.Range("BF78").Value = 0
$wsMeasures.Range("BG78").Value = 0
$wsMeasures.Range("BH78").Value = 0
$wsMeasures.Range("BI78").Value = 0

# New ingredient / recipe rows on NutritionalData sheet (rows 208-211)
$wsNutri.Range("A208").Value = 'Good&Gather broccoli'
$wsNutri.Range("B208").Formula = "=30*3.5"
$wsNutri.Range("C208").Formula = "=0"
$wsNutri.Range("D208").Value = 0
$wsNutri.Range("E208").Formula = "=1*3.5"
$wsNutri.Range("F208").Formula = "=4*3.5"
$wsNutri.Range("G208").Formula = "=2*3.5"
$wsNutri.Range("H208").Formula = "=20*3.5"

$wsNutri.Range("A209").Value = 'penne red fennel barilla/1 pkg broccoli frozen Good & Gather brand/1 pkg beyond meat/1 yellow bell pepper/2 tbs olive oil/4 tbs sourcream Winco brand 1 pot makes about 5 bowls'
$wsNutri.Range("B209").Formula = "=SUM(B208,B204,B6*2,B61*5,B37*4.5,B50*2)"
$wsNutri.Range("C209").Formula = "=SUM(C208,C204,C6*2,C61*5,C37*4.5,C50*2)"
$wsNutri.Range("D209").Formula = "=SUM(D208,D204,D6*2,D61*5,D37*4.5,D50*2)"
$wsNutri.Range("E209").Formula = "=SUM(E208,E204,E6*2,E61*5,E37*4.5,E50*2)"
$wsNutri.Range("F209").Formula = "=SUM(F208,F204,F6*2,F61*5,F37*4.5,F50*2)"
$wsNutri.Range("G209").Formula = "=SUM(G208,G204,G6*2,G61*5,G37*4.5,G50*2)"
$wsNutri.Range("H209").Formula = "=SUM(H208,H204,H6*2,H61*5,H37*4.5,H50*2)"

$wsNutri.Range("A210").Value = 'bowl of 4/1/2021 pasta recipe'
$wsNutri.Range("B210").Formula = "=B209/5"
$wsNutri.Range("C210").Formula = "=C209/5"
$wsNutri.Range("D210").Formula = "=D209/5"
$wsNutri.Range("E210").Formula = "=E209/5"
$wsNutri.Range("F210").Formula = "=F209/5"
$wsNutri.Range("G210").Formula = "=G209/5"
$wsNutri.Range("H210").Formula = "=H209/5"

$wsNutri.Range("A211").Value = 'Starbucks Nitro cold brew with sweet cream https://www.starbucks.com/menu/product/2122237/iced?parent=%2Fdrinks%2Fcold-coffees%2Fnitro-cold-brews'
$wsNutri.Range("B211").Value = 70
$wsNutri.Range("C211").Value = 5
$wsNutri.Range("D211").Value = 3.5
$wsNutri.Range("E211").Value = 1
$wsNutri.Range("F211").Value = 4
$wsNutri.Range("G211").Value = 0
$wsNutri.Range("H211").Value = 20
